$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header C1 from "URL name" to "Row"
$ws.Range("C1").Value = "Row"

# Fill column C (rows 2-84) with a =ROW() formula
for ($r = 2; $r -le 84; $r++) {
    $ws.Cells.Item($r, 3).Formula = "=ROW()"
}

# Update D27 formula to also reference B$19
$ws.Range("D27").Formula = '="["&B6&","&B$19&"]"'

# Restore the view: frozen pane top-left cell and active selection
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Range("D28").Select()
